$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6, shifting the existing
# rows 6-10 (with their data) down to rows 7-11.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly price record.
$ws.Cells.Item(6, 1).Value  = 1
$ws.Cells.Item(6, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value  = 44539
$ws.Cells.Item(6, 5).Value  = 15
$ws.Cells.Item(6, 6).Value  = 100112017
$ws.Cells.Item(6, 7).Value  = "Ramas de apio"
$ws.Cells.Item(6, 8).Value  = "Americana (o)"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(6, 11).Value = 6500
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 6750
$ws.Cells.Item(6, 14).Value = "$/atado 7 kilos"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 6750
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
